$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'58.668.95"
$ws.Range('E2').Value = "'  -3.93%  "
$ws.Range('D3').Value = "'2.627.73"
$ws.Range('E3').Value = "'  -3.43%  "
$ws.Range('E4').Value = "'  -0.08%  "
$ws.Range('D5').Value = "'525.00"
$ws.Range('E5').Value = "'  -0.84%  "
$ws.Range('D6').Value = "'143.47"
$ws.Range('E6').Value = "'  -3.30%  "
$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = "'  +0.24%  "
$ws.Range('D8').Value = "'0.571"
$ws.Range('E8').Value = "'  -1.50%  "
$ws.Range('D9').Value = "'6.60"
$ws.Range('E9').Value = "'  -7.09%  "
$ws.Range('E10').Value = "'  -2.95%  "
$ws.Range('E11').Value = "'  -1.89%  "
$ws.Range('E12').Value = "'  +0.94%  "
$ws.Range('D13').Value = "'3.089.55"
$ws.Range('D14').Value = "'58.606.83"
$ws.Range('E14').Value = "'  -4.00%  "
$ws.Range('D15').Value = "'20.88"
$ws.Range('E15').Value = "'  -3.17%  "
$ws.Range('B16').Value = "'ShibaInu"
$ws.Range('C16').Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range('D16').Value = "'0.0000136"
$ws.Range('E16').Value = "'  -1.99%  "
$ws.Range('B17').Value = "'WrappedEther"
$ws.Range('C17').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D17').Value = "'2.639.20"
$ws.Range('E17').Value = "'  -3.60%  "
$ws.Range('B18').Value = "'Polkadot"
$ws.Range('C18').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('D18').Value = "'4.44"
$ws.Range('E18').Value = "'  -1.56%  "
$ws.Range('B19').Value = "'BitcoinCash"
$ws.Range('C19').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('D19').Value = "'338.78"
$ws.Range('E19').Value = "'  -2.24%  "
$ws.Range('E20').Value = "'  -1.13%  "
$ws.Range('E21').Value = "'  -1.54%  "
$ws.Range('E22').Value = "'  +0.04%  "
$ws.Range('D23').Value = "'65.38"
$ws.Range('E23').Value = "'  +3.06%  "
$ws.Range('E24').Value = "'  -0.79%  "
$ws.Range('E25').Value = "'  -2.95%  "
$ws.Range('E26').Value = "'  +0.52%  "
$ws.Range('D27').Value = "'7.19"
$ws.Range('E27').Value = "'  -2.05%  "
$ws.Range('E28').Value = "'  -4.18%  "
$ws.Range('E29').Value = "'  -3.67%  "
$ws.Range('D30').Value = "'0.998"
$ws.Range('E30').Value = "'  +0.05%  "
$ws.Range('E31').Value = "'  -0.40%  "
$ws.Range('D32').Value = "'18.88"
$ws.Range('E32').Value = "'  -1.21%  "
$ws.Range('D33').Value = "'150.00"
$ws.Range('E33').Value = "'  -0.28%  "
$ws.Range('D34').Value = "'4.13"
$ws.Range('E34').Value = "'  -3.09%  "
$ws.Range('E35').Value = "'  -3.46%  "
$ws.Range('D36').Value = "'0.900"
$ws.Range('E36').Value = "'  -2.55%  "
$ws.Range('E37').Value = "'  -4.89%  "
$ws.Range('D38').Value = "'36.39"
$ws.Range('E38').Value = "'  -3.54%  "
$ws.Range('E39').Value = "'  -6.77%  "
$ws.Range('E40').Value = "'  -1.33%  "
$ws.Range('D41').Value = "'0.997"
$ws.Range('E41').Value = "'  +0.22%  "
$ws.Range('D42').Value = "'0.603"
$ws.Range('E42').Value = "'  -4.18%  "
$ws.Range('D43').Value = "'0.0973"
$ws.Range('E43').Value = "'  -1.29%  "
$ws.Range('D44').Value = "'270.89"
$ws.Range('E44').Value = "'  -3.52%  "
$ws.Range('E45').Value = "'  +1.25%  "
$ws.Range('B46').Value = "'Hedera"
$ws.Range('C46').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('D46').Value = "'0.0534"
$ws.Range('E46').Value = "'  -1.44%  "
$ws.Range('B47').Value = "'EnergySwap"
$ws.Range('C47').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('D47').Value = "'19.12"
$ws.Range('E47').Value = "'  -5.70%  "
$ws.Range('D48').Value = "'2.040.38"
$ws.Range('E48').Value = "'  -2.95%  "
$ws.Range('E49').Value = "'  -2.00%  "
$ws.Range('D50').Value = "'4.61"
$ws.Range('E50').Value = "'  -7.63%  "
$ws.Range('D51').Value = "'18.53"
$ws.Range('E51').Value = "'  -4.75%  "
